# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# leve-profit tables (one table per crafting job sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 104167050
$ws.Range("I98").Value = 113636590
$ws.Range("J98").Value = 2006
$ws.Range("K98").Value = 113636590
$ws.Range("L98").Value = 2006
$ws.Range("M98").Value = -113635092
$ws.Range("N98").Value = -5002

$ws.Range("H122").Value = 104167050
$ws.Range("I122").Value = 113636590
$ws.Range("J122").Value = 2006
$ws.Range("K122").Value = 340909770
$ws.Range("L122").Value = 6018
$ws.Range("M122").Value = -340907320
$ws.Range("N122").Value = -10918

$ws.Range("H132").Value = 29416378
$ws.Range("I132").Value = 24594364
$ws.Range("J132").Value = 71436780
$ws.Range("K132").Value = 73783092
$ws.Range("L132").Value = 214310340
$ws.Range("M132").Value = -73780562
$ws.Range("N132").Value = -214315400

$ws.Range("H138").Value = 3864.9404
$ws.Range("I138").Value = 2594.6
$ws.Range("J138").Value = 4261.922
$ws.Range("K138").Value = 7783.799999999999
$ws.Range("L138").Value = 12785.766
$ws.Range("M138").Value = -2643.799999999999
$ws.Range("N138").Value = -23065.766

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 911.44446
$ws.Range("I2").Value = 881.875
$ws.Range("J2").Value = 954.4545000000001
$ws.Range("K2").Value = 881.875
$ws.Range("L2").Value = 954.4545000000001
$ws.Range("M2").Value = -768.875
$ws.Range("N2").Value = -1180.4545

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H97").Value = 1480.5625
$ws.Range("I97").Value = 1288.625
$ws.Range("J97").Value = 1672.5
$ws.Range("K97").Value = 1288.625
$ws.Range("L97").Value = 1672.5
$ws.Range("M97").Value = -792.625
$ws.Range("N97").Value = -2664.5

$ws.Range("H116").Value = 911.44446
$ws.Range("I116").Value = 881.875
$ws.Range("J116").Value = 954.4545000000001
$ws.Range("K116").Value = 881.875
$ws.Range("L116").Value = 954.4545000000001
$ws.Range("M116").Value = 1412.125
$ws.Range("N116").Value = -5542.4545

$ws.Range("H132").Value = 21742652
$ws.Range("I132").Value = 30304194
$ws.Range("J132").Value = 9504.538
$ws.Range("K132").Value = 90912582
$ws.Range("L132").Value = 28513.614
$ws.Range("M132").Value = -90910052
$ws.Range("N132").Value = -33573.614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 911.44446
$ws.Range("I3").Value = 881.875
$ws.Range("J3").Value = 954.4545000000001
$ws.Range("K3").Value = 881.875
$ws.Range("L3").Value = 954.4545000000001
$ws.Range("M3").Value = -767.875
$ws.Range("N3").Value = -1182.4545

$ws.Range("H86").Value = 1663623.4
$ws.Range("I86").Value = 3283.2
$ws.Range("J86").Value = 2586034.8
$ws.Range("K86").Value = 3283.2
$ws.Range("L86").Value = 2586034.8
$ws.Range("M86").Value = -2160.2
$ws.Range("N86").Value = -2588280.8

$ws.Range("H89").Value = 1663623.4
$ws.Range("I89").Value = 3283.2
$ws.Range("J89").Value = 2586034.8
$ws.Range("K89").Value = 16416
$ws.Range("L89").Value = 12930174
$ws.Range("M89").Value = -10800
$ws.Range("N89").Value = -12941406

$ws.Range("H107").Value = 55556236
$ws.Range("I107").Value = 100000570
$ws.Range("J107").Value = 818.25
$ws.Range("K107").Value = 100000570
$ws.Range("L107").Value = 818.25
$ws.Range("M107").Value = -99998650
$ws.Range("N107").Value = -4658.25

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H134").Value = 3713123
$ws.Range("I134").Value = 7119.857
$ws.Range("J134").Value = 12360463
$ws.Range("K134").Value = 21359.571
$ws.Range("L134").Value = 37081389
$ws.Range("M134").Value = -18824.571
$ws.Range("N134").Value = -37086459

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9528106
$ws.Range("I132").Value = 870.24
$ws.Range("J132").Value = 33346194
$ws.Range("K132").Value = 2610.72
$ws.Range("L132").Value = 100038582
$ws.Range("M132").Value = -80.72000000000025
$ws.Range("N132").Value = -100043642

$ws.Range("H134").Value = 1871.4839
$ws.Range("I134").Value = 2104.8333
$ws.Range("J134").Value = 1071.4286
$ws.Range("K134").Value = 6314.499899999999
$ws.Range("L134").Value = 3214.2858
$ws.Range("M134").Value = -3779.499899999999
$ws.Range("N134").Value = -8284.2858

$ws.Range("H137").Value = 33816.668
$ws.Range("J137").Value = 33816.668
$ws.Range("L137").Value = 33816.668
$ws.Range("N137").Value = -44016.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 41675896
$ws.Range("I122").Value = 125000420
$ws.Range("J122").Value = 13637.375
$ws.Range("K122").Value = 1125003780
$ws.Range("L122").Value = 122736.375
$ws.Range("M122").Value = -1125001330
$ws.Range("N122").Value = -127636.375

$ws.Range("H131").Value = 729.16
$ws.Range("I131").Value = 440.9091
$ws.Range("J131").Value = 764.7865
$ws.Range("K131").Value = 1322.7273
$ws.Range("L131").Value = 2294.3595
$ws.Range("M131").Value = 3717.2727
$ws.Range("N131").Value = -12374.3595

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13894886
$ws.Range("I122").Value = 20840340
$ws.Range("J122").Value = 3975.25
$ws.Range("K122").Value = 62521020
$ws.Range("L122").Value = 11925.75
$ws.Range("M122").Value = -62518570
$ws.Range("N122").Value = -16825.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 82419290
$ws.Range("I16").Value = 7937611
$ws.Range("J16").Value = 250003060
$ws.Range("K16").Value = 7937611
$ws.Range("L16").Value = 250003060
$ws.Range("M16").Value = -7937441
$ws.Range("N16").Value = -250003400

$ws.Range("H22").Value = 1267331.8
$ws.Range("I22").Value = 3164779.2
$ws.Range("J22").Value = 2366.6667
$ws.Range("K22").Value = 3164779.2
$ws.Range("L22").Value = 2366.6667
$ws.Range("M22").Value = -3164484.2
$ws.Range("N22").Value = -2956.6667

$ws.Range("H27").Value = 1267331.8
$ws.Range("I27").Value = 3164779.2
$ws.Range("J27").Value = 2366.6667
$ws.Range("K27").Value = 3164779.2
$ws.Range("L27").Value = 2366.6667
$ws.Range("M27").Value = -3164672.2
$ws.Range("N27").Value = -2580.6667

$ws.Range("H100").Value = 4420.4
$ws.Range("J100").Value = 3900.5
$ws.Range("L100").Value = 3900.5
$ws.Range("N100").Value = -4982.5

$ws.Range("H125").Value = 40715
$ws.Range("J125").Value = 40715
$ws.Range("L125").Value = 40715
$ws.Range("N125").Value = -50555

$ws.Range("H127").Value = 40548.5
$ws.Range("J127").Value = 40548.5
$ws.Range("L127").Value = 40548.5
$ws.Range("N127").Value = -50468.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2299.5908
$ws.Range("I122").Value = 1564.1765
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 4692.529500000001
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -2242.529500000001
$ws.Range("N122").Value = -19300

$ws.Range("H132").Value = 16149070
$ws.Range("I132").Value = 23835550
$ws.Range("J132").Value = 7461.9
$ws.Range("K132").Value = 71506650
$ws.Range("L132").Value = 22385.7
$ws.Range("M132").Value = -71504120
$ws.Range("N132").Value = -27445.7
